$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 45: KKR vs PBKS - enter per-team scores for row 57
$ws.Range("E57").Value = 60
$ws.Range("H57").Value = 50
$ws.Range("K57").Value = 30
$ws.Range("N57").Value = 0
$ws.Range("Q57").Value = 80
$ws.Range("T57").Value = 40
$ws.Range("W57").Value = 20
$ws.Range("Z57").Value = 100
$ws.Range("AC57").Value = 70
